$wb = $excel.ActiveWorkbook

# --- Update "Yearly" sheet ---
# June 2017: Taxable Account (L8) and 401K (M8) dividend amounts were revised.
# The dependent formula cells (O8, L15, M15, O15) recalc automatically.
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("L8").Value = 152.4
$wsYearly.Range("M8").Value = 58.75

# --- Update "All Time" sheet ---
# Totals for 2017 (row 8) and the grand totals (row 46) pull from Yearly via
# formulas and recalc automatically once the Yearly inputs above change.
$wsAllTime = $wb.Worksheets.Item("All Time")

# Scroll the All Time view down before picking the new selection, matching
# the author's final on-screen view/selection when the edit was made.
$wsAllTime.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("M53").Select()

# Leave the Yearly sheet with its own updated selection.
$wsYearly.Activate()
$wsYearly.Range("M8").Select()

# Re-activate "All Time" so it matches the workbook's last-active sheet.
$wsAllTime.Activate()

$wb.Save()
